$wb = $excel.ActiveWorkbook

# --- Update the text note on "Hoja1"!A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.04 = 24111.29 pesos`n✅ 24111.29 pesos = 6.03 = 970.39 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate values on "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 165.69
$ws2.Range("O10").Value = 3995
$ws2.Range("N12").Value = 3998.99
$ws2.Range("O12").Value = 160.945
